$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Bmp2"
$ws.Range("C2").Value = "Bmpr1a"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 9.163165666666666
$ws.Range("H2").Value = 27.489497
$ws.Range("I2").Value = 0.2800251397703982
$ws.Range("J2").Value = 0.2800251397703982
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.750436666666666
$ws.Range("N2").Value = 14.25131
$ws.Range("O2").Value = 0.07745299862590357
$ws.Range("P2").Value = 0.07745299862590359
$ws.Range("Q2").Value = 43.52903816567444
$ws.Range("R2").Value = 391.76134349107
$ws.Range("S2").Value = 0.02168878676585511
$ws.Range("T2").Value = 0.02168878676585511

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Bmp2"
$ws.Range("C3").Value = "Bmpr1a"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 9.163165666666666
$ws.Range("H3").Value = 27.489497
$ws.Range("I3").Value = 0.2800251397703982
$ws.Range("J3").Value = 0.2800251397703982
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 34.05277366666667
$ws.Range("N3").Value = 102.158321
$ws.Range("O3").Value = 0.5552098927072401
$ws.Range("P3").Value = 0.5552098927072401
$ws.Range("Q3").Value = 312.0312065171707
$ws.Range("R3").Value = 2808.280858654537
$ws.Range("S3").Value = 0.1554727278072527
$ws.Range("T3").Value = 0.1554727278072527

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Bmp2"
$ws.Range("C4").Value = "Bmpr1a"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 9.163165666666666
$ws.Range("H4").Value = 27.489497
$ws.Range("I4").Value = 0.2800251397703982
$ws.Range("J4").Value = 0.2800251397703982
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 22.52994333333334
$ws.Range("N4").Value = 67.58983
$ws.Range("O4").Value = 0.3673371086668564
$ws.Range("P4").Value = 0.3673371086668564
$ws.Range("Q4").Value = 206.4456032239456
$ws.Range("R4").Value = 1858.01042901551
$ws.Range("S4").Value = 0.1028636251972904
$ws.Range("T4").Value = 0.1028636251972904

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Bmp2"
$ws.Range("C5").Value = "Bmpr1a"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 12.06704733333333
$ws.Range("H5").Value = 36.201142
$ws.Range("I5").Value = 0.3687673822623249
$ws.Range("J5").Value = 0.3687673822623249
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.750436666666666
$ws.Range("N5").Value = 14.25131
$ws.Range("O5").Value = 0.07745299862590357
$ws.Range("P5").Value = 0.07745299862590359
$ws.Range("Q5").Value = 57.32374411066888
$ws.Range("R5").Value = 515.91369699602
$ws.Range("S5").Value = 0.02856213955164191
$ws.Range("T5").Value = 0.02856213955164191

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Bmp2"
$ws.Range("C6").Value = "Bmpr1a"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 12.06704733333333
$ws.Range("H6").Value = 36.201142
$ws.Range("I6").Value = 0.3687673822623249
$ws.Range("J6").Value = 0.3687673822623249
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 34.05277366666667
$ws.Range("N6").Value = 102.158321
$ws.Range("O6").Value = 0.5552098927072401
$ws.Range("P6").Value = 0.5552098927072401
$ws.Range("Q6").Value = 410.9164316669535
$ws.Range("R6").Value = 3698.247885002582
$ws.Range("S6").Value = 0.2047432987397952
$ws.Range("T6").Value = 0.2047432987397952

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Bmp2"
$ws.Range("C7").Value = "Bmpr1a"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 12.06704733333333
$ws.Range("H7").Value = 36.201142
$ws.Range("I7").Value = 0.3687673822623249
$ws.Range("J7").Value = 0.3687673822623249
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 22.52994333333334
$ws.Range("N7").Value = 67.58983
$ws.Range("O7").Value = 0.3673371086668564
$ws.Range("P7").Value = 0.3673371086668564
$ws.Range("Q7").Value = 271.8698926206511
$ws.Range("R7").Value = 2446.82903358586
$ws.Range("S7").Value = 0.1354619439708878
$ws.Range("T7").Value = 0.1354619439708878

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Bmp2"
$ws.Range("C8").Value = "Bmpr1a"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 11.49244066666667
$ws.Range("H8").Value = 34.477322
$ws.Range("I8").Value = 0.3512074779672769
$ws.Range("J8").Value = 0.351207477967277
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 4.750436666666666
$ws.Range("N8").Value = 14.25131
$ws.Range("O8").Value = 0.07745299862590357
$ws.Range("P8").Value = 0.07745299862590359
$ws.Range("Q8").Value = 54.59411153242444
$ws.Range("R8").Value = 491.34700379182
$ws.Range("S8").Value = 0.02720207230840656
$ws.Range("T8").Value = 0.02720207230840657

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Bmp2"
$ws.Range("C9").Value = "Bmpr1a"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 11.49244066666667
$ws.Range("H9").Value = 34.477322
$ws.Range("I9").Value = 0.3512074779672769
$ws.Range("J9").Value = 0.351207477967277
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 34.05277366666667
$ws.Range("N9").Value = 102.158321
$ws.Range("O9").Value = 0.5552098927072401
$ws.Range("P9").Value = 0.5552098927072401
$ws.Range("Q9").Value = 391.3494808995958
$ws.Range("R9").Value = 3522.145328096362
$ws.Range("S9").Value = 0.1949938661601922
$ws.Range("T9").Value = 0.1949938661601922

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Bmp2"
$ws.Range("C10").Value = "Bmpr1a"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 11.49244066666667
$ws.Range("H10").Value = 34.477322
$ws.Range("I10").Value = 0.3512074779672769
$ws.Range("J10").Value = 0.351207477967277
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 22.52994333333334
$ws.Range("N10").Value = 67.58983
$ws.Range("O10").Value = 0.3673371086668564
$ws.Range("P10").Value = 0.3673371086668564
$ws.Range("Q10").Value = 258.9240369816956
$ws.Range("R10").Value = 2330.31633283526
$ws.Range("S10").Value = 0.1290115394986782
$ws.Range("T10").Value = 0.1290115394986782
